$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry (matching original inlineStr/text cells) by prefixing with a literal
# single quote, Excel's standard "treat as text" marker. This prevents values that
# look numeric (e.g. "64.196.94", "0.0000118") from being auto-converted to numbers.

$ws.Range("D2").Value = "'64.196.94"
$ws.Range("E2").Value = "'  -2.69%  "
$ws.Range("D3").Value = "'3.171.81"
$ws.Range("E3").Value = "'  -7.79%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'559.93"
$ws.Range("E5").Value = "'  -4.08%  "
$ws.Range("D6").Value = "'171.82"
$ws.Range("E6").Value = "'  -1.01%  "
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "'  -0.36%  "
$ws.Range("D9").Value = "'3.170.00"
$ws.Range("E9").Value = "'  -7.85%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "'  -5.83%  "
$ws.Range("D11").Value = "'6.63"
$ws.Range("E11").Value = "'  -4.34%  "
$ws.Range("D12").Value = "'0.396"
$ws.Range("E12").Value = "'  -3.18%  "
$ws.Range("D13").Value = "'3.716.63"
$ws.Range("E13").Value = "'  -7.96%  "
$ws.Range("E14").Value = "'  +0.40%  "
$ws.Range("D15").Value = "'27.53"
$ws.Range("E15").Value = "'  -4.03%  "
$ws.Range("D16").Value = "'64.140.87"
$ws.Range("E16").Value = "'  -2.83%  "
$ws.Range("E17").Value = "'  -4.91%  "
$ws.Range("D18").Value = "'3.168.89"
$ws.Range("E18").Value = "'  -8.17%  "
$ws.Range("D19").Value = "'5.67"
$ws.Range("E19").Value = "'  -4.45%  "
$ws.Range("D20").Value = "'13.04"
$ws.Range("E20").Value = "'  -5.80%  "
$ws.Range("D21").Value = "'352.55"
$ws.Range("E21").Value = "'  -4.47%  "
$ws.Range("D22").Value = "'7.17"
$ws.Range("E22").Value = "'  -6.35%  "
$ws.Range("E23").Value = "'  -0.19%  "
$ws.Range("D24").Value = "'69.20"
$ws.Range("E24").Value = "'  -4.37%  "
$ws.Range("B25").Value = "'PEPE"
$ws.Range("C25").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000118"
$ws.Range("E25").Value = "'  -2.98%  "
$ws.Range("B26").Value = "'Polygon"
$ws.Range("C26").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "'0.502"
$ws.Range("E26").Value = "'  -5.98%  "
$ws.Range("D27").Value = "'9.45"
$ws.Range("E27").Value = "'  -3.00%  "
$ws.Range("E28").Value = "'  -0.58%  "
$ws.Range("E29").Value = "'  -0.20%  "
$ws.Range("B30").Value = "'NEARProtocol"
$ws.Range("C30").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'5.62"
$ws.Range("E30").Value = "'  -2.01%  "
$ws.Range("B31").Value = "'USDe"
$ws.Range("C31").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  +0.04%  "
$ws.Range("D32").Value = "'1.89"
$ws.Range("E32").Value = "'  -4.69%  "
$ws.Range("D33").Value = "'22.08"
$ws.Range("E33").Value = "'  -6.58%  "
$ws.Range("D34").Value = "'6.61"
$ws.Range("E34").Value = "'  -5.60%  "
$ws.Range("D35").Value = "'1.20"
$ws.Range("E35").Value = "'  -6.54%  "
$ws.Range("D36").Value = "'157.04"
$ws.Range("E36").Value = "'  -2.17%  "
$ws.Range("E37").Value = "'  -5.87%  "
$ws.Range("D38").Value = "'26.03"
$ws.Range("E38").Value = "'  -9.47%  "
$ws.Range("E39").Value = "'  -9.55%  "
$ws.Range("D40").Value = "'2.54"
$ws.Range("E40").Value = "'  -2.68%  "
$ws.Range("E41").Value = "'  -4.93%  "
$ws.Range("D42").Value = "'2.650.54"
$ws.Range("E42").Value = "'  -4.23%  "
$ws.Range("D43").Value = "'6.04"
$ws.Range("E43").Value = "'  -6.61%  "
$ws.Range("D44").Value = "'4.15"
$ws.Range("E44").Value = "'  -6.72%  "
$ws.Range("D45").Value = "'0.0652"
$ws.Range("E45").Value = "'  -4.16%  "
$ws.Range("D46").Value = "'38.78"
$ws.Range("E46").Value = "'  -3.64%  "
$ws.Range("D47").Value = "'326.76"
$ws.Range("E47").Value = "'  +0.50%  "
$ws.Range("D48").Value = "'23.64"
$ws.Range("E48").Value = "'  -3.20%  "
$ws.Range("E49").Value = "'  -6.76%  "
$ws.Range("E50").Value = "'  -0.94%  "
$ws.Range("E51").Value = "'  -0.06%  "
